$d = $word.ActiveDocument
$d.Content.Find.Execute("Hello, World", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
